# Append a new data row (row 31) to each of the 4 worksheets.
# Row 31 mirrors the structure of the existing rows: a timestamp in
# column A (formatted as YYYY-MM-DD HH:MM:SS) plus the raw/decoded
# hex-byte fields in columns B-I.

$wb = $excel.ActiveWorkbook

$dateValue = [double]"45817.43892361111"

function Add-DataRow {
    param(
        $ws,
        $rowNum,
        $colB,
        $colC,
        $colD,
        $colE,
        $colF,
        $colG,
        $colH,
        $colI
    )

    $ws.Cells.Item($rowNum, 1).Value = $dateValue
    $ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowNum, 2).Value = $colB
    $ws.Cells.Item($rowNum, 3).Value = $colC
    $ws.Cells.Item($rowNum, 4).Value = $colD
    $ws.Cells.Item($rowNum, 5).Value = $colE

    $ws.Cells.Item($rowNum, 6).Value = $colF
    $ws.Cells.Item($rowNum, 7).Value = [double]$colG
    $ws.Cells.Item($rowNum, 8).Value = $colH
    $ws.Cells.Item($rowNum, 9).Value = $colI
}

# Sheet 1: DE_LFT_#1
$ws1 = $wb.Worksheets.Item("DE_LFT_#1")
Add-DataRow $ws1 31 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x78" "0x14" 380 "7.598631275147109e+23" 376 14

# Sheet 2: DE_LFT_#2
$ws2 = $wb.Worksheets.Item("DE_LFT_#2")
Add-DataRow $ws2 31 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x78" "0xe" 380 "5.68432987514711e+23" 376 14

# Sheet 3: DE_PLT_#1
$ws3 = $wb.Worksheets.Item("DE_PLT_#1")
Add-DataRow $ws3 31 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x82" "0x7" 130 "5.68631262647114e+23" 129 7

# Sheet 4: DE_PLT_#2
$ws4 = $wb.Worksheets.Item("DE_PLT_#2")
Add-DataRow $ws4 31 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x81" "0x3" 130 "9.85046333984776e+23" 129 3
